$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") rows 2-13 from date serial 45212 to 45221
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 3).Value = 45221
}
